$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(3).Delete()
